# Adds "2022-Q4" fund-holding data to the workbook.
#
# Before:  Sheet1 = "总计", Sheet2 = "2022-Q2" (fund holdings table)
# After:   Sheet1 = "总计", Sheet2 = "2022-Q4" (new fund holdings table),
#          Sheet3 = "2022-Q2" (the original fund holdings table, unchanged)
#
# Strategy:
#   1. Duplicate the existing "2022-Q2" sheet - the duplicate keeps the old
#      data/formatting exactly and ends up named "2022-Q2" again.
#   2. Rename the original sheet (now redundant) to "2022-Q4" and overwrite
#      its contents with the new Q4 fund table.
#   3. Insert a new "2022-Q4" summary row on the "总计" sheet, pushing the
#      pre-existing "2022-Q2" row down.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: duplicate "2022-Q2" (sheet index 2) so its data + styling is
# preserved verbatim in a new sheet.
# ------------------------------------------------------------------
$origQ2 = $wb.Worksheets.Item(2)
$origQ2.Copy($null, $origQ2)

# Free up the "2022-Q2" name on the original sheet so the duplicate can
# take it.
$origQ2.Name = "2022-Q4"

$newQ2 = $wb.Worksheets.Item(3)
$newQ2.Name = "2022-Q2"

# ------------------------------------------------------------------
# Step 2: overwrite the (renamed) "2022-Q4" sheet with the new fund table.
# ------------------------------------------------------------------
$q4 = $origQ2
$q4.Cells.Clear()

# -- Header row, styled like the "总计" header (bold / centered / bordered). --
$totalSheet = $wb.Worksheets.Item(1)
$totalSheet.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$q4.Range("A2:A10").PasteSpecial(-4122) # same style down column A

$headerCols = @("B","C","D","E","F","G","H")
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $addr = $headerCols[$i] + "1"
    $q4.Range($addr).Value = $headers[$i]
}

# -- Data rows. Columns D,E,F,G and B hold numeric-looking text, so force
# text formatting before assigning, then strip the format override back off
# (matching the source file's un-styled data cells). --
$textCols = @("B","D","E","F","G")
foreach ($col in $textCols) {
    $q4.Range($col + "2:" + $col + "10").NumberFormat = "@"
}

$rows = @(
    @("001230","鹏华医药科技股票","16.95","94.13","3.79","0.6424",10),
    @("012093","鹏华创新升级混合A","5.64","75.32","4.02","0.2267",8),
    @("014313","鹏华创新增长一年持有期混合A","3.47","60.62","3.96","0.1374",5),
    @("005108","圆信永丰双利优选定期开放灵活配置混合","0.73","88.26","5.31","0.0388",4),
    @("001965","圆信永丰兴源灵活配置混合A","0.58","84.58","4.96","0.0288",5),
    @("006274","圆信永丰医药健康混合","0.24","85.59","5.13","0.0123",5),
    @("012094","鹏华创新升级混合C","0.17","75.32","4.02","0.0068",8),
    @("014314","鹏华创新增长一年持有期混合C","0.15","60.62","3.96","0.0059",5),
    @("001966","圆信永丰兴源灵活配置混合C","0.09","84.58","4.96","0.0045",5)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]
    $q4.Range("A" + $rowNum).Value = $r
    $q4.Range("B" + $rowNum).Value = $data[0]
    $q4.Range("C" + $rowNum).Value = $data[1]
    $q4.Range("D" + $rowNum).Value = $data[2]
    $q4.Range("E" + $rowNum).Value = $data[3]
    $q4.Range("F" + $rowNum).Value = $data[4]
    $q4.Range("G" + $rowNum).Value = $data[5]
    $q4.Range("H" + $rowNum).Value = $data[6]
}

foreach ($col in $textCols) {
    $q4.Range($col + "2:" + $col + "10").ClearFormats()
}

# ------------------------------------------------------------------
# Step 3: update the "总计" sheet - new row for 2022-Q4, old 2022-Q2 row
# pushed down to row 3.
# ------------------------------------------------------------------
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)  # xlPasteFormats -> s="2"

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 1.1
